# Generate Report for Handoff
# Renames the handed-off file from the old GUID to a new GUID, refreshes the
# handoff timestamps / xliff hash, and clears out the per-locale "handback"
# bookkeeping columns (Latest Target File / Latest Handback File / Latest
# Handback DateTime) now that a fresh handoff round has started.

$wb = $excel.ActiveWorkbook

$oldGuid = "60fbb8d4-661e-4a21-9688-2ae0bcc04d4b"
$newGuid = "eff7a7ef-e4dc-43d9-82bf-37a005ba4a35"
$newHash = "ee344444aeece6df9bc22de073153dc0f76a7ea8"

$newFileName       = "$newGuid.md"
$newPathAndName    = "e2e\$newGuid.md"
$newHoGenDate      = "2016-08-16 00:54:45"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName

$wsOverview.Hyperlinks.Delete()
$ovLink = $wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f90fd891fe89046f7c575c6a21930debdef12886/e2e/$newGuid.md")
$ovLink.TextToDisplay = $newPathAndName

$wsOverview.Range("G2").Value = $newHoGenDate

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-16 00:54:40"

# Target/handback bookkeeping is reset for the new handoff round.
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Only the Source File Name hyperlink (A2) survives; the Latest Target File
# hyperlink (old I2) is gone now that the column is empty.
$wsZhCn.Hyperlinks.Delete()
$zhLink = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f90fd891fe89046f7c575c6a21930debdef12886/e2e/$newGuid.md")
$zhLink.TextToDisplay = $newFileName

$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = $newHoGenDate

# Target/handback bookkeeping is reset for the new handoff round.
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Hyperlinks.Delete()
$deLink = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f90fd891fe89046f7c575c6a21930debdef12886/e2e/$newGuid.md")
$deLink.TextToDisplay = $newFileName

$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
